# Auto-generated Excel COM-interop script applying the scheduled-runner price/profit
# refresh across every Ragnarok_Profits crafting-job sheet (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# Each write sets currentAveragePrice/NQ/HQ (H:L) and the derived LeveProfit (M:N) cells
# to their newly recalculated values for the affected leve rows.
$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(28, 8).Value = 3301.4783
$ws.Cells.Item(28, 9).Value = 633.4211
$ws.Cells.Item(28, 10).Value = 15974.75
$ws.Cells.Item(28, 11).Value = 633.4211
$ws.Cells.Item(28, 12).Value = 15974.75
$ws.Cells.Item(28, 13).Value = -148.4211
$ws.Cells.Item(28, 14).Value = -16944.75
$ws.Cells.Item(62, 8).Value = 4808.647
$ws.Cells.Item(62, 9).Value = 2925.3
$ws.Cells.Item(62, 10).Value = 7499.143
$ws.Cells.Item(62, 11).Value = 2925.3
$ws.Cells.Item(62, 12).Value = 7499.143
$ws.Cells.Item(62, 13).Value = -2301.3
$ws.Cells.Item(62, 14).Value = -8747.143
$ws.Cells.Item(65, 8).Value = 4808.647
$ws.Cells.Item(65, 9).Value = 2925.3
$ws.Cells.Item(65, 10).Value = 7499.143
$ws.Cells.Item(65, 11).Value = 14626.5
$ws.Cells.Item(65, 12).Value = 37495.715
$ws.Cells.Item(65, 13).Value = -11506.5
$ws.Cells.Item(65, 14).Value = -43735.715
$ws.Cells.Item(94, 8).Value = 2212.3572
$ws.Cells.Item(94, 9).Value = 2212.3572
$ws.Cells.Item(94, 11).Value = 2212.3572
$ws.Cells.Item(94, 13).Value = -1761.3572
$ws.Cells.Item(111, 8).Value = 2215.9092
$ws.Cells.Item(111, 9).Value = 1708.4445
$ws.Cells.Item(111, 11).Value = 5125.333500000001
$ws.Cells.Item(111, 13).Value = -2058.333500000001
$ws.Cells.Item(112, 8).Value = 3017.1943
$ws.Cells.Item(112, 10).Value = 3200.6667
$ws.Cells.Item(112, 12).Value = 9602.000100000001
$ws.Cells.Item(112, 14).Value = -11818.0001
$ws.Cells.Item(132, 8).Value = 1283.9412
$ws.Cells.Item(132, 9).Value = 918.5833
$ws.Cells.Item(132, 11).Value = 2755.7499
$ws.Cells.Item(132, 13).Value = -225.7498999999998
$ws.Cells.Item(138, 8).Value = 6269.4614
$ws.Cells.Item(138, 9).Value = 3391.6428
$ws.Cells.Item(138, 10).Value = 9626.916999999999
$ws.Cells.Item(138, 11).Value = 10174.9284
$ws.Cells.Item(138, 12).Value = 28880.751
$ws.Cells.Item(138, 13).Value = -5034.928400000001
$ws.Cells.Item(138, 14).Value = -39160.751

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 4867.67
$ws.Cells.Item(32, 9).Value = 4707.084
$ws.Cells.Item(32, 11).Value = 4707.084
$ws.Cells.Item(32, 13).Value = -4420.084
$ws.Cells.Item(45, 8).Value = 2372.611
$ws.Cells.Item(45, 9).Value = 1694.5
$ws.Cells.Item(45, 10).Value = 4135.7
$ws.Cells.Item(45, 11).Value = 1694.5
$ws.Cells.Item(45, 12).Value = 4135.7
$ws.Cells.Item(45, 13).Value = -1317.5
$ws.Cells.Item(45, 14).Value = -4889.7
$ws.Cells.Item(74, 8).Value = 1317493
$ws.Cells.Item(74, 9).Value = 1472088.6
$ws.Cells.Item(74, 10).Value = 3430
$ws.Cells.Item(74, 11).Value = 1472088.6
$ws.Cells.Item(74, 12).Value = 3430
$ws.Cells.Item(74, 13).Value = -1471214.6
$ws.Cells.Item(74, 14).Value = -5178
$ws.Cells.Item(77, 8).Value = 1317493
$ws.Cells.Item(77, 9).Value = 1472088.6
$ws.Cells.Item(77, 10).Value = 3430
$ws.Cells.Item(77, 11).Value = 7360443
$ws.Cells.Item(77, 12).Value = 17150
$ws.Cells.Item(77, 13).Value = -7356075
$ws.Cells.Item(77, 14).Value = -25886
$ws.Cells.Item(110, 8).Value = 7531.154
$ws.Cells.Item(110, 10).Value = 5599.2
$ws.Cells.Item(110, 12).Value = 5599.2
$ws.Cells.Item(110, 14).Value = -9689.200000000001
$ws.Cells.Item(122, 8).Value = 4618.6113
$ws.Cells.Item(122, 9).Value = 2939.077
$ws.Cells.Item(122, 11).Value = 8817.231
$ws.Cells.Item(122, 13).Value = -6367.231
$ws.Cells.Item(132, 8).Value = 6670391.5
$ws.Cells.Item(132, 9).Value = 3881.1667
$ws.Cells.Item(132, 11).Value = 11643.5001
$ws.Cells.Item(132, 13).Value = -9113.500100000001

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(99, 8).Value = 1572.6364
$ws.Cells.Item(99, 9).Value = 1385.1578
$ws.Cells.Item(99, 11).Value = 1385.1578
$ws.Cells.Item(99, 13).Value = 112.8422
$ws.Cells.Item(105, 8).Value = 1148026.9
$ws.Cells.Item(105, 9).Value = 1761457
$ws.Cells.Item(105, 11).Value = 1761457
$ws.Cells.Item(105, 13).Value = -1759710

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(6, 8).Value = 253150.25
$ws.Cells.Item(6, 9).Value = 999999
$ws.Cells.Item(6, 10).Value = 4200.6665
$ws.Cells.Item(6, 11).Value = 999999
$ws.Cells.Item(6, 12).Value = 4200.6665
$ws.Cells.Item(6, 13).Value = -999886
$ws.Cells.Item(6, 14).Value = -4426.6665
$ws.Cells.Item(7, 8).Value = 29.434782
$ws.Cells.Item(7, 10).Value = 40.2
$ws.Cells.Item(7, 12).Value = 40.2
$ws.Cells.Item(7, 14).Value = -266.2
$ws.Cells.Item(31, 8).Value = 11104154
$ws.Cells.Item(31, 9).Value = 22731312
$ws.Cells.Item(31, 10).Value = 219153.31
$ws.Cells.Item(31, 11).Value = 22731312
$ws.Cells.Item(31, 12).Value = 219153.31
$ws.Cells.Item(31, 13).Value = -22731017
$ws.Cells.Item(31, 14).Value = -219743.31
$ws.Cells.Item(34, 8).Value = 11104154
$ws.Cells.Item(34, 9).Value = 22731312
$ws.Cells.Item(34, 10).Value = 219153.31
$ws.Cells.Item(34, 11).Value = 22731312
$ws.Cells.Item(34, 12).Value = 219153.31
$ws.Cells.Item(34, 13).Value = -22731110
$ws.Cells.Item(34, 14).Value = -219557.31
$ws.Cells.Item(38, 8).Value = 14532.5
$ws.Cells.Item(38, 9).Value = 24065
$ws.Cells.Item(38, 10).Value = 5000
$ws.Cells.Item(38, 11).Value = 24065
$ws.Cells.Item(38, 12).Value = 5000
$ws.Cells.Item(38, 13).Value = -23655
$ws.Cells.Item(38, 14).Value = -5820
$ws.Cells.Item(134, 8).Value = 1332.0714
$ws.Cells.Item(134, 9).Value = 1383.7894
$ws.Cells.Item(134, 11).Value = 4151.3682
$ws.Cells.Item(134, 13).Value = -1616.3682

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 536.6111
$ws.Cells.Item(5, 9).Value = 395.7143
$ws.Cells.Item(5, 10).Value = 626.2727
$ws.Cells.Item(5, 11).Value = 1187.1429
$ws.Cells.Item(5, 12).Value = 1878.8181
$ws.Cells.Item(5, 13).Value = -1075.1429
$ws.Cells.Item(5, 14).Value = -2102.8181
$ws.Cells.Item(32, 8).Value = 1686213.1
$ws.Cells.Item(32, 10).Value = 1823455.8
$ws.Cells.Item(32, 12).Value = 5470367.4
$ws.Cells.Item(32, 14).Value = -5470933.4
$ws.Cells.Item(68, 8).Value = 1932.8334
$ws.Cells.Item(68, 9).Value = 1149.5
$ws.Cells.Item(68, 10).Value = 3499.5
$ws.Cells.Item(68, 11).Value = 3448.5
$ws.Cells.Item(68, 12).Value = 10498.5
$ws.Cells.Item(68, 13).Value = -2637.5
$ws.Cells.Item(68, 14).Value = -12120.5
$ws.Cells.Item(71, 8).Value = 1932.8334
$ws.Cells.Item(71, 9).Value = 1149.5
$ws.Cells.Item(71, 10).Value = 3499.5
$ws.Cells.Item(71, 11).Value = 10345.5
$ws.Cells.Item(71, 12).Value = 31495.5
$ws.Cells.Item(71, 13).Value = -6289.5
$ws.Cells.Item(71, 14).Value = -39607.5
$ws.Cells.Item(107, 8).Value = 2938461.2
$ws.Cells.Item(107, 9).Value = 1404
$ws.Cells.Item(107, 10).Value = 4793444.5
$ws.Cells.Item(107, 11).Value = 4212
$ws.Cells.Item(107, 12).Value = 14380333.5
$ws.Cells.Item(107, 13).Value = -2292
$ws.Cells.Item(107, 14).Value = -14384173.5
$ws.Cells.Item(108, 8).Value = 13466
$ws.Cells.Item(108, 9).Value = 8499.25
$ws.Cells.Item(108, 11).Value = 25497.75
$ws.Cells.Item(108, 13).Value = -22617.75
$ws.Cells.Item(131, 8).Value = 5248.2964
$ws.Cells.Item(131, 10).Value = 5618.8
$ws.Cells.Item(131, 12).Value = 16856.4
$ws.Cells.Item(131, 14).Value = -26936.4
$ws.Cells.Item(134, 8).Value = 5551.864
$ws.Cells.Item(134, 9).Value = 2343.9285
$ws.Cells.Item(134, 10).Value = 11165.75
$ws.Cells.Item(134, 11).Value = 7031.7855
$ws.Cells.Item(134, 12).Value = 33497.25
$ws.Cells.Item(134, 13).Value = -1961.7855
$ws.Cells.Item(134, 14).Value = -43637.25
$ws.Cells.Item(135, 8).Value = 536.6111
$ws.Cells.Item(135, 9).Value = 395.7143
$ws.Cells.Item(135, 10).Value = 626.2727
$ws.Cells.Item(135, 11).Value = 3561.4287
$ws.Cells.Item(135, 12).Value = 5636.454299999999
$ws.Cells.Item(135, 13).Value = -1026.4287
$ws.Cells.Item(135, 14).Value = -10706.4543
$ws.Cells.Item(140, 8).Value = 3294.6785
$ws.Cells.Item(140, 9).Value = 1082.7727
$ws.Cells.Item(140, 11).Value = 3248.3181
$ws.Cells.Item(140, 13).Value = 1931.6819

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 8479.862999999999
$ws.Cells.Item(70, 9).Value = 5947.0835
$ws.Cells.Item(70, 11).Value = 5947.0835
$ws.Cells.Item(70, 13).Value = -5677.0835
$ws.Cells.Item(73, 8).Value = 8479.862999999999
$ws.Cells.Item(73, 9).Value = 5947.0835
$ws.Cells.Item(73, 11).Value = 5947.0835
$ws.Cells.Item(73, 13).Value = -5011.0835
$ws.Cells.Item(80, 8).Value = 3229.8333
$ws.Cells.Item(80, 9).Value = 2718.625
$ws.Cells.Item(80, 10).Value = 4252.25
$ws.Cells.Item(80, 11).Value = 2718.625
$ws.Cells.Item(80, 12).Value = 4252.25
$ws.Cells.Item(80, 13).Value = -1720.625
$ws.Cells.Item(80, 14).Value = -6248.25
$ws.Cells.Item(83, 8).Value = 3229.8333
$ws.Cells.Item(83, 9).Value = 2718.625
$ws.Cells.Item(83, 10).Value = 4252.25
$ws.Cells.Item(83, 11).Value = 13593.125
$ws.Cells.Item(83, 12).Value = 21261.25
$ws.Cells.Item(83, 13).Value = -8601.125
$ws.Cells.Item(83, 14).Value = -31245.25
$ws.Cells.Item(106, 8).Value = 64749
$ws.Cells.Item(106, 10).Value = 64749
$ws.Cells.Item(106, 12).Value = 64749
$ws.Cells.Item(106, 14).Value = -67273

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(98, 8).Value = 99973
$ws.Cells.Item(98, 10).Value = 99973
$ws.Cells.Item(98, 12).Value = 99973
$ws.Cells.Item(98, 14).Value = -105963
$ws.Cells.Item(100, 10).Value = 25028314
$ws.Cells.Item(100, 12).Value = 25028314
$ws.Cells.Item(100, 14).Value = -25029396
$ws.Cells.Item(122, 8).Value = 3828.8333
$ws.Cells.Item(122, 9).Value = 3369.8438
$ws.Cells.Item(122, 10).Value = 7500.75
$ws.Cells.Item(122, 11).Value = 10109.5314
$ws.Cells.Item(122, 12).Value = 22502.25
$ws.Cells.Item(122, 13).Value = -7659.5314
$ws.Cells.Item(122, 14).Value = -27402.25
$ws.Cells.Item(136, 8).Value = 2657.4
$ws.Cells.Item(136, 9).Value = 1378.3125
$ws.Cells.Item(136, 11).Value = 4134.9375
$ws.Cells.Item(136, 13).Value = -1584.9375

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 5728.0625
$ws.Cells.Item(81, 10).Value = 5621.4287
$ws.Cells.Item(81, 12).Value = 11242.8574
$ws.Cells.Item(81, 14).Value = -13364.8574
$ws.Cells.Item(84, 8).Value = 5728.0625
$ws.Cells.Item(84, 10).Value = 5621.4287
$ws.Cells.Item(84, 12).Value = 56214.287
$ws.Cells.Item(84, 14).Value = -66822.28700000001
$ws.Cells.Item(96, 8).Value = 7132
$ws.Cells.Item(96, 9).Value = 7123.75
$ws.Cells.Item(96, 10).Value = 7143
$ws.Cells.Item(96, 11).Value = 7123.75
$ws.Cells.Item(96, 12).Value = 7143
$ws.Cells.Item(96, 13).Value = -5750.75
$ws.Cells.Item(96, 14).Value = -9889
$ws.Cells.Item(136, 8).Value = 359289.1
$ws.Cells.Item(136, 9).Value = 1477
$ws.Cells.Item(136, 10).Value = 772149.25
$ws.Cells.Item(136, 11).Value = 4431
$ws.Cells.Item(136, 12).Value = 2316447.75
$ws.Cells.Item(136, 13).Value = -1881
$ws.Cells.Item(136, 14).Value = -2321547.75
